$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '258.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.12%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.02'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-3.28%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.883'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-8.73%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05962'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.687'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8756'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.66%'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9622'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5.53%'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1414'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.49%'

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07182'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.10%'

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03137'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.37%'

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09235'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.08%'

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001545'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.18%'

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'One'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006048'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.04%'

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006004'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.66%'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.485'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.44%'

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.69%'

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3145'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.68%'

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03598'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4.25%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.529'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.07%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04224'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.71%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1380'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.11%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001222'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.58%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004519'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.29%'

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.07%'

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001493'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '2.64%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03837'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.40%'

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1104'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.10%'

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004015'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-30.25%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002309'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.61%'

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '5.59%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005489'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '4.20%'

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.04%'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1091'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '9.14%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002156'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1.74%'

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'
